# Auto-generated Excel COM-interop script to apply scheduled-runner value updates
# to the Chocobo_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 407002
$ws.Range("I116").Value = 836041.56
$ws.Range("J116").Value = 10965.462
$ws.Range("K116").Value = 836041.56
$ws.Range("L116").Value = 10965.462
$ws.Range("M116").Value = -832599.56
$ws.Range("N116").Value = -17849.462

$ws.Range("H124").Value = 46535.715
$ws.Range("J124").Value = 46535.715
$ws.Range("L124").Value = 46535.715
$ws.Range("N124").Value = -56355.715

$ws.Range("H132").Value = 445856.1
$ws.Range("I132").Value = 244817.52
$ws.Range("J132").Value = 2506501.5
$ws.Range("K132").Value = 734452.5599999999
$ws.Range("L132").Value = 7519504.5
$ws.Range("M132").Value = -731922.5599999999
$ws.Range("N132").Value = -7524564.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 12672.667
$ws.Range("I28").Value = 3749.5
$ws.Range("J28").Value = 30519
$ws.Range("K28").Value = 3749.5
$ws.Range("L28").Value = 30519
$ws.Range("M28").Value = -3557.5
$ws.Range("N28").Value = -30903

$ws.Range("H61").Value = 1000.6222
$ws.Range("I61").Value = 783.60974
$ws.Range("J61").Value = 3225
$ws.Range("K61").Value = 783.60974
$ws.Range("L61").Value = 3225
$ws.Range("M61").Value = -571.60974
$ws.Range("N61").Value = -3649

$ws.Range("H99").Value = 12672.667
$ws.Range("I99").Value = 3749.5
$ws.Range("J99").Value = 30519
$ws.Range("K99").Value = 3749.5
$ws.Range("L99").Value = 30519
$ws.Range("M99").Value = -754.5
$ws.Range("N99").Value = -36509

$ws.Range("H102").Value = 1284.5938
$ws.Range("I102").Value = 1214.6522
$ws.Range("J102").Value = 1463.3334
$ws.Range("K102").Value = 1214.6522
$ws.Range("L102").Value = 1463.3334
$ws.Range("M102").Value = 407.3478
$ws.Range("N102").Value = -4707.3334

$ws.Range("H136").Value = 1000.6222
$ws.Range("I136").Value = 783.60974
$ws.Range("J136").Value = 3225
$ws.Range("K136").Value = 2350.82922
$ws.Range("L136").Value = 9675
$ws.Range("M136").Value = 199.1707799999999
$ws.Range("N136").Value = -14775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1688.4814
$ws.Range("I105").Value = 1650.2916
$ws.Range("K105").Value = 1650.2916
$ws.Range("M105").Value = 96.70839999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11630171
$ws.Range("I31").Value = 1279.5518
$ws.Range("J31").Value = 35718588
$ws.Range("K31").Value = 1279.5518
$ws.Range("L31").Value = 35718588
$ws.Range("M31").Value = -984.5518
$ws.Range("N31").Value = -35719178

$ws.Range("H34").Value = 11630171
$ws.Range("I34").Value = 1279.5518
$ws.Range("J34").Value = 35718588
$ws.Range("K34").Value = 1279.5518
$ws.Range("L34").Value = 35718588
$ws.Range("M34").Value = -1077.5518
$ws.Range("N34").Value = -35718992

$ws.Range("H58").Value = 1480.4788
$ws.Range("I58").Value = 1297.5366
$ws.Range("J58").Value = 2730.5833
$ws.Range("K58").Value = 1297.5366
$ws.Range("L58").Value = 2730.5833
$ws.Range("M58").Value = -1094.5366
$ws.Range("N58").Value = -3136.5833

$ws.Range("H68").Value = 56260.145
$ws.Range("J68").Value = 56260.145
$ws.Range("L68").Value = 56260.145
$ws.Range("N68").Value = -57758.145

$ws.Range("H71").Value = 56260.145
$ws.Range("J71").Value = 56260.145
$ws.Range("L71").Value = 168780.435
$ws.Range("N71").Value = -176268.435

$ws.Range("H99").Value = 10532195
$ws.Range("I99").Value = 25003964
$ws.Range("J99").Value = 7272.727
$ws.Range("K99").Value = 25003964
$ws.Range("L99").Value = 7272.727
$ws.Range("M99").Value = -25002466
$ws.Range("N99").Value = -10268.727

$ws.Range("H126").Value = 10532195
$ws.Range("I126").Value = 25003964
$ws.Range("J126").Value = 7272.727
$ws.Range("K126").Value = 75011892
$ws.Range("L126").Value = 21818.181
$ws.Range("M126").Value = -75009422
$ws.Range("N126").Value = -26758.181

$ws.Range("H134").Value = 1419.0209
$ws.Range("I134").Value = 649.13794
$ws.Range("J134").Value = 2594.1052
$ws.Range("K134").Value = 1947.41382
$ws.Range("L134").Value = 7782.3156
$ws.Range("M134").Value = 587.5861800000002
$ws.Range("N134").Value = -12852.3156

$ws.Range("H135").Value = 39580
$ws.Range("J135").Value = 39580
$ws.Range("L135").Value = 39580
$ws.Range("N135").Value = -49720

$ws.Range("H136").Value = 1480.4788
$ws.Range("I136").Value = 1297.5366
$ws.Range("J136").Value = 2730.5833
$ws.Range("K136").Value = 3892.6098
$ws.Range("L136").Value = 8191.749899999999
$ws.Range("M136").Value = -1342.6098
$ws.Range("N136").Value = -13291.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 4582.353
$ws.Range("I112").Value = 3866.6667
$ws.Range("J112").Value = 4735.7144
$ws.Range("K112").Value = 11600.0001
$ws.Range("L112").Value = 14207.1432
$ws.Range("M112").Value = -10492.0001
$ws.Range("N112").Value = -16423.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 70000
$ws.Range("J25").Value = 70000
$ws.Range("L25").Value = 70000
$ws.Range("N25").Value = -71058

$ws.Range("H38").Value = 29999.5
$ws.Range("J38").Value = 29999.5
$ws.Range("L38").Value = 29999.5
$ws.Range("N38").Value = -30925.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 20049
$ws.Range("J39").Value = 20049
$ws.Range("L39").Value = 20049
$ws.Range("N39").Value = -20875

$ws.Range("H42").Value = 31683
$ws.Range("J42").Value = 35024.5
$ws.Range("L42").Value = 35024.5
$ws.Range("N42").Value = -35780.5

$ws.Range("H43").Value = 16442.5
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 15256.667
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 15256.667
$ws.Range("M43").Value = -19851
$ws.Range("N43").Value = -15554.667

$ws.Range("H46").Value = 66466.664
$ws.Range("J46").Value = 66466.664
$ws.Range("L46").Value = 66466.664
$ws.Range("N46").Value = -66928.664

$ws.Range("H112").Value = 31140
$ws.Range("J112").Value = 31140
$ws.Range("L112").Value = 31140
$ws.Range("N112").Value = -34094

$ws.Range("H122").Value = 3175.8918
$ws.Range("I122").Value = 2126.0454
$ws.Range("J122").Value = 4715.6665
$ws.Range("K122").Value = 6378.1362
$ws.Range("L122").Value = 14146.9995
$ws.Range("M122").Value = -3928.1362
$ws.Range("N122").Value = -19046.9995

$ws.Range("H134").Value = 66466.664
$ws.Range("J134").Value = 66466.664
$ws.Range("L134").Value = 199399.992
$ws.Range("N134").Value = -204469.992
